$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B
$ws.Range("B2").Value = 12.56424632954657
$ws.Range("B3").Value = 12.3361462690672
$ws.Range("B4").Value = 12.19740464054708
$ws.Range("B5").Value = 12.14127811622863
$ws.Range("B6").Value = 12.13198561452287
$ws.Range("B7").Value = 12.19664592307438
$ws.Range("B8").Value = 12.48537387139891
$ws.Range("B9").Value = 13.05842839606041
$ws.Range("B10").Value = 13.47897443982496
$ws.Range("B11").Value = 13.66924673285522
$ws.Range("B12").Value = 13.7410745728849
$ws.Range("B13").Value = 13.72561628915044
$ws.Range("B14").Value = 13.67516096224707
$ws.Range("B15").Value = 13.64422425859104
$ws.Range("B16").Value = 13.46651234184087
$ws.Range("B17").Value = 13.35717185231178
$ws.Range("B18").Value = 13.29418884084202
$ws.Range("B19").Value = 13.27285020675971
$ws.Range("B20").Value = 13.36882154117985
$ws.Range("B21").Value = 13.68998757897338
$ws.Range("B22").Value = 13.89854940854388
$ws.Range("B23").Value = 13.7873824042152
$ws.Range("B24").Value = 13.36355509012418
$ws.Range("B25").Value = 12.90315667981841

# Column C
$ws.Range("C2").Value = 8.693692551266818
$ws.Range("C3").Value = 8.524058822954572
$ws.Range("C4").Value = 8.420999557194019
$ws.Range("C5").Value = 8.379340802398403
$ws.Range("C6").Value = 8.372445700541107
$ws.Range("C7").Value = 8.420436277796465
$ws.Range("C8").Value = 8.635012955184946
$ws.Range("C9").Value = 9.061757596461273
$ws.Range("C10").Value = 9.375337237370482
$ws.Range("C11").Value = 9.51728260110373
$ws.Range("C12").Value = 9.570875780018968
$ws.Range("C13").Value = 9.559341457399425
$ws.Range("C14").Value = 9.52169523725644
$ws.Range("C15").Value = 9.498613556294695
$ws.Range("C16").Value = 9.366041653181675
$ws.Range("C17").Value = 9.28449120754984
$ws.Range("C18").Value = 9.237522572177264
$ws.Range("C19").Value = 9.221610776004679
$ws.Range("C20").Value = 9.293179333691223
$ws.Range("C21").Value = 9.532757586667588
$ws.Range("C22").Value = 9.688387832403347
$ws.Range("C23").Value = 9.605429926913263
$ws.Range("C24").Value = 9.289251689705004
$ws.Range("C25").Value = 8.946050405476502

# Column E
$ws.Range("E2").Value = 18.43564304162264
$ws.Range("E3").Value = 18.18317526234162
$ws.Range("E4").Value = 18.03124594427196
$ws.Range("E5").Value = 17.97018449487464
$ws.Range("E6").Value = 17.96009883639077
$ws.Range("E7").Value = 18.03041890509877
$ws.Range("E8").Value = 18.34799677249703
$ws.Range("E9").Value = 18.99191208037906
$ws.Range("E10").Value = 19.4733655858584
$ws.Range("E11").Value = 19.69325388233212
$ws.Range("E12").Value = 19.7765671211532
$ws.Range("E13").Value = 19.7586233050833
$ws.Range("E14").Value = 19.70010763122038
$ws.Range("E15").Value = 19.66426880489048
$ws.Range("E16").Value = 19.45900611943344
$ws.Range("E17").Value = 19.33325073198023
$ws.Range("E18").Value = 19.26100633885014
$ws.Range("E19").Value = 19.23656297491127
$ws.Range("E20").Value = 19.34662920566807
$ws.Range("E21").Value = 19.71729446774796
$ws.Range("E22").Value = 19.95978029651085
$ws.Range("E23").Value = 19.83036522239995
$ws.Range("E24").Value = 19.34058062331105
$ws.Range("E25").Value = 18.81591911387293

# Column F
$ws.Range("F2").Value = 45.94816114895949
$ws.Range("F3").Value = 45.72393092305035
$ws.Range("F4").Value = 45.59837716250062
$ws.Range("F5").Value = 45.55029707374767
$ws.Range("F6").Value = 45.54250061958123
$ws.Range("F7").Value = 45.59771620737944
$ws.Range("F8").Value = 45.86835273299541
$ws.Range("F9").Value = 46.49350297544158
$ws.Range("F10").Value = 47.00769885669116
$ws.Range("F11").Value = 47.25288528199575
$ws.Range("F12").Value = 47.34729165682041
$ws.Range("F13").Value = 47.32689121083746
$ws.Range("F14").Value = 47.26062126056502
$ws.Range("F15").Value = 47.22023021601998
$ws.Range("F16").Value = 46.99189703549657
$ws.Range("F17").Value = 46.85466712916562
$ws.Range("F18").Value = 46.77680140565371
$ws.Range("F19").Value = 46.7506222060923
$ws.Range("F20").Value = 46.8691656714379
$ws.Range("F21").Value = 47.28004453352717
$ws.Range("F22").Value = 47.55763725731681
$ws.Range("F23").Value = 47.40867309705101
$ws.Range("F24").Value = 46.86260766853218
$ws.Range("F25").Value = 46.31453233828372

# Column G
$ws.Range("G2").Value = 3.691439664340859
$ws.Range("G3").Value = 3.694167971353137
$ws.Range("G4").Value = 3.695929424181401
$ws.Range("G5").Value = 3.696668998545062
$ws.Range("G6").Value = 3.696793121300999
$ws.Range("G7").Value = 3.695939310100486
$ws.Range("G8").Value = 3.692362527477654
$ws.Range("G9").Value = 3.68602939686989
$ws.Range("G10").Value = 3.681786630126445
$ws.Range("G11").Value = 3.679944501647058
$ws.Range("G12").Value = 3.679259499145707
$ws.Range("G13").Value = 3.679406468713183
$ws.Range("G14").Value = 3.679887894538763
$ws.Range("G15").Value = 3.680184416809273
$ws.Range("G16").Value = 3.681908780690902
$ws.Range("G17").Value = 3.682989090842589
$ws.Range("G18").Value = 3.683618737103643
$ws.Range("G19").Value = 3.683833348892028
$ws.Range("G20").Value = 3.682873233512783
$ws.Range("G21").Value = 3.679746147488886
$ws.Range("G22").Value = 3.67777565949907
$ws.Range("G23").Value = 3.678820668018366
$ws.Range("G24").Value = 3.682925585935767
$ws.Range("G25").Value = 3.687670289297589

# Column I
$ws.Range("I2").Value = 26.47732755246626
$ws.Range("I3").Value = 26.54212489604378
$ws.Range("I4").Value = 26.58722995984778
$ws.Range("I5").Value = 26.60694395221487
$ws.Range("I6").Value = 26.6102978471185
$ws.Range("I7").Value = 26.58749043693114
$ws.Range("I8").Value = 26.49856315643127
$ws.Range("I9").Value = 26.36656793476753
$ws.Range("I10").Value = 26.2956794908289
$ws.Range("I11").Value = 26.26914453108809
$ws.Range("I12").Value = 26.25992158514182
$ws.Range("I13").Value = 26.26187114919889
$ws.Range("I14").Value = 26.26836919313683
$ws.Range("I15").Value = 26.27245701623109
$ws.Range("I16").Value = 26.29752889450831
$ws.Range("I17").Value = 26.31437557915756
$ws.Range("I18").Value = 26.3246028238972
$ws.Range("I19").Value = 26.32815778599458
$ws.Range("I20").Value = 26.31252656557149
$ws.Range("I21").Value = 26.26643813330126
$ws.Range("I22").Value = 26.24112858165141
$ws.Range("I23").Value = 26.25419527396868
$ws.Range("I24").Value = 26.31336081689748
$ws.Range("I25").Value = 26.39771181029196

# Column J
$ws.Range("J2").Value = 9.217611969698449
$ws.Range("J3").Value = 9.24246373974313
$ws.Range("J4").Value = 9.258878284394713
$ws.Range("J5").Value = 9.265858009048564
$ws.Range("J6").Value = 9.267034545930269
$ws.Range("J7").Value = 9.258971238359129
$ws.Range("J8").Value = 9.225941113472015
$ws.Range("J9").Value = 9.170332167729477
$ws.Range("J10").Value = 9.135055667853415
$ws.Range("J11").Value = 9.120217737591357
$ws.Range("J12").Value = 9.11477286338156
$ws.Range("J13").Value = 9.115937778907396
$ws.Range("J14").Value = 9.119766298716742
$ws.Range("J15").Value = 9.122134028300676
$ws.Range("J16").Value = 9.136049700061488
$ws.Range("J17").Value = 9.14489630022976
$ws.Range("J18").Value = 9.150098476607575
$ws.Range("J19").Value = 9.151879397424775
$ws.Range("J20").Value = 9.143942781932568
$ws.Range("J21").Value = 9.118637049259652
$ws.Range("J22").Value = 9.103112079856905
$ws.Range("J23").Value = 9.111305284794565
$ws.Range("J24").Value = 9.144373505756173
$ws.Range("J25").Value = 9.184395607917351

# Column K
$ws.Range("K2").Value = 13.28618600433387
$ws.Range("K3").Value = 13.13701660721303
$ws.Range("K4").Value = 13.04805500529264
$ws.Range("K5").Value = 13.01250446067726
$ws.Range("K6").Value = 13.00664488821105
$ws.Range("K7").Value = 13.04757266275455
$ws.Range("K8").Value = 13.23423173916297
$ws.Range("K9").Value = 13.61926749934121
$ws.Range("K10").Value = 13.91117634545188
$ws.Range("K11").Value = 14.04537929524874
$ws.Range("K12").Value = 14.09635524928911
$ws.Range("K13").Value = 14.08537045886048
$ws.Range("K14").Value = 14.04957025366421
$ws.Range("K15").Value = 14.02766055477623
$ws.Range("K16").Value = 13.90243026794538
$ws.Range("K17").Value = 13.82593376520746
$ws.Range("K18").Value = 13.78207070496367
$ws.Range("K19").Value = 13.76724416836627
$ws.Range("K20").Value = 13.83406324435973
$ws.Range("K21").Value = 14.06008177327869
$ws.Range("K22").Value = 14.20868712218546
$ws.Range("K23").Value = 14.12930763288988
$ws.Range("K24").Value = 13.83038754250644
$ws.Range("K25").Value = 13.51333306587298

# Column N
$ws.Range("N2").Value = 20.62553880119735
$ws.Range("N3").Value = 20.69348743264106
$ws.Range("N4").Value = 20.73718695628315
$ws.Range("N5").Value = 20.7554935964667
$ws.Range("N6").Value = 20.75856355601467
$ws.Range("N7").Value = 20.73743182521583
$ws.Range("N8").Value = 20.64855748470955
$ws.Range("N9").Value = 20.48992247762674
$ws.Range("N10").Value = 20.38283584928125
$ws.Range("N11").Value = 20.33615789488375
$ws.Range("N12").Value = 20.3187738375806
$ws.Range("N13").Value = 20.32250483913673
$ws.Range("N14").Value = 20.33472185384001
$ws.Range("N15").Value = 20.34224311062124
$ws.Range("N16").Value = 20.38592725747855
$ws.Range("N17").Value = 20.41324693048482
$ws.Range("N18").Value = 20.42915220410735
$ws.Range("N19").Value = 20.43457041379678
$ws.Range("N20").Value = 20.41031886918781
$ws.Range("N21").Value = 20.33112550660833
$ws.Range("N22").Value = 20.28106888984945
$ws.Range("N23").Value = 20.30762972972952
$ws.Range("N24").Value = 20.41164202619847
$ws.Range("N25").Value = 20.53117027865397
